# live_trading_results.xlsx update
# Trade #60 closed at 2026-02-16 21:34:02 - leadlag DOWN +0.000%
#
# This script:
#  1) Updates the Summary sheet aggregate stats (OVERALL / leadlag / momentum rows)
#  2) Marks trades #32-#34 (leadlag sheet rows 28-30) as CLOSED with exit data
#  3) Appends a brand-new OPEN trade #60 to the leadlag sheet (row 49)
#  4) Marks trades #35-#38 (momentum sheet rows 7-10) as CLOSED with exit data
#  5) Appends the now-closed trades (#32-#38, across leadlag+momentum) to the
#     "All Trades" sheet as rows 33-39
#  6) Updates the Comparison sheet stats for leadlag / momentum

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    param($ws, [string]$addr, [string]$text)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
}

function Set-EmptyTextCell {
    param($ws, [string]$addr)
    # A leading apostrophe forces Excel to treat the (otherwise empty) entry
    # as literal text, producing a real empty-string text cell rather than a
    # blank/uninitialised one.
    $ws.Range($addr).Value = "'"
}

function Set-NumCell {
    param($ws, [string]$addr, $num)
    $ws.Range($addr).Value = $num
}

# ---------------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

Set-NumCell  $wsSummary "C2" 38
Set-TextCell $wsSummary "D2" "63.2%"
Set-TextCell $wsSummary "E2" "+7.5493%"
Set-TextCell $wsSummary "F2" "+0.1987%"

Set-NumCell  $wsSummary "C3" 47
Set-TextCell $wsSummary "D3" "36.2%"
Set-TextCell $wsSummary "E3" "+4.6292%"
Set-TextCell $wsSummary "F3" "+0.0985%"

Set-NumCell  $wsSummary "C4" 12
Set-TextCell $wsSummary "D4" "58.3%"
Set-TextCell $wsSummary "E4" "+2.9201%"
Set-TextCell $wsSummary "F4" "+0.2433%"

# ---------------------------------------------------------------------------
# 2) leadlag sheet - close out trades #32, #33, #34 (rows 28, 29, 30)
# ---------------------------------------------------------------------------
$wsLeadlag = $wb.Worksheets.Item("leadlag")

Set-NumCell  $wsLeadlag "G28" 69085.281617
Set-TextCell $wsLeadlag "H28" "CLOSED"
Set-NumCell  $wsLeadlag "I28" -0.2841
Set-NumCell  $wsLeadlag "J28" -2.84
Set-TextCell $wsLeadlag "M28" "time_exit_5min"
Set-NumCell  $wsLeadlag "N28" 5

Set-NumCell  $wsLeadlag "G29" 68930.30664900001
Set-TextCell $wsLeadlag "H29" "CLOSED"
Set-NumCell  $wsLeadlag "I29" -0.1174
Set-NumCell  $wsLeadlag "J29" -1.17
Set-TextCell $wsLeadlag "M29" "time_exit_5min"
Set-NumCell  $wsLeadlag "N29" 5

Set-NumCell  $wsLeadlag "G30" 68703.52499799999
Set-TextCell $wsLeadlag "H30" "CLOSED"
Set-NumCell  $wsLeadlag "I30" 0.1552
Set-NumCell  $wsLeadlag "J30" 1.55
Set-TextCell $wsLeadlag "M30" "time_exit_5min"
Set-NumCell  $wsLeadlag "N30" 5

# New trade #60 - freshly opened, appended as row 49
Set-NumCell     $wsLeadlag "A49" 60
Set-TextCell    $wsLeadlag "B49" "2026-02-16"
Set-TextCell    $wsLeadlag "C49" "21:34:02"
Set-TextCell    $wsLeadlag "D49" "leadlag"
Set-TextCell    $wsLeadlag "E49" "DOWN"
Set-NumCell     $wsLeadlag "F49" 68769.89999999999
Set-EmptyTextCell $wsLeadlag "G49"
Set-TextCell    $wsLeadlag "H49" "OPEN"
Set-NumCell     $wsLeadlag "I49" 0
Set-NumCell     $wsLeadlag "J49" 0
Set-NumCell     $wsLeadlag "K49" 0.75
Set-TextCell    $wsLeadlag "L49" "Coinbase leading with -0.086% move"
Set-EmptyTextCell $wsLeadlag "M49"
Set-NumCell     $wsLeadlag "N49" 0

# ---------------------------------------------------------------------------
# 3) momentum sheet - close out trades #35, #36, #37, #38 (rows 7, 8, 9, 10)
# ---------------------------------------------------------------------------
$wsMomentum = $wb.Worksheets.Item("momentum")

Set-NumCell  $wsMomentum "G7" 69046.8011
Set-TextCell $wsMomentum "H7" "CLOSED"
Set-NumCell  $wsMomentum "I7" -0.4771
Set-NumCell  $wsMomentum "J7" -4.77
Set-TextCell $wsMomentum "M7" "time_exit_5min"
Set-NumCell  $wsMomentum "N7" 5

Set-NumCell  $wsMomentum "G8" 68269.16643899999
Set-TextCell $wsMomentum "H8" "CLOSED"
Set-NumCell  $wsMomentum "I8" 0.6915
Set-NumCell  $wsMomentum "J8" 6.92
Set-TextCell $wsMomentum "M8" "time_exit_5min"
Set-NumCell  $wsMomentum "N8" 5

Set-NumCell  $wsMomentum "G9" 68423.353225
Set-TextCell $wsMomentum "H9" "CLOSED"
Set-NumCell  $wsMomentum "I9" 0.232
Set-NumCell  $wsMomentum "J9" 2.32
Set-TextCell $wsMomentum "M9" "time_exit_5min"
Set-NumCell  $wsMomentum "N9" 5

Set-NumCell  $wsMomentum "G10" 68082.94107
Set-TextCell $wsMomentum "H10" "CLOSED"
Set-NumCell  $wsMomentum "I10" 0.7307
Set-NumCell  $wsMomentum "J10" 7.31
Set-TextCell $wsMomentum "M10" "time_exit_5min"
Set-NumCell  $wsMomentum "N10" 5

# ---------------------------------------------------------------------------
# 4) All Trades sheet - append the now-closed trades as rows 33-39
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

# Row 33 <- leadlag trade #32 (leadlag row 28)
Set-NumCell     $wsAll "A33" 32
Set-TextCell    $wsAll "B33" "2026-02-16"
Set-TextCell    $wsAll "C33" "21:28:22"
Set-TextCell    $wsAll "D33" "leadlag"
Set-TextCell    $wsAll "E33" "DOWN"
Set-NumCell     $wsAll "F33" 68889.565
Set-NumCell     $wsAll "G33" 69085.281617
Set-TextCell    $wsAll "H33" "CLOSED"
Set-NumCell     $wsAll "I33" -0.2841
Set-NumCell     $wsAll "J33" -2.84
Set-NumCell     $wsAll "K33" 0.7052
Set-TextCell    $wsAll "L33" "Coinbase leading with -0.071% move"
Set-TextCell    $wsAll "M33" "time_exit_5min"
Set-NumCell     $wsAll "N33" 5

# Row 34 <- leadlag trade #33 (leadlag row 29)
Set-NumCell     $wsAll "A34" 33
Set-TextCell    $wsAll "B34" "2026-02-16"
Set-TextCell    $wsAll "C34" "21:28:28"
Set-TextCell    $wsAll "D34" "leadlag"
Set-TextCell    $wsAll "E34" "DOWN"
Set-NumCell     $wsAll "F34" 68849.49000000001
Set-NumCell     $wsAll "G34" 68930.30664900001
Set-TextCell    $wsAll "H34" "CLOSED"
Set-NumCell     $wsAll "I34" -0.1174
Set-NumCell     $wsAll "J34" -1.17
Set-NumCell     $wsAll "K34" 0.75
Set-TextCell    $wsAll "L34" "Binance leading with -0.160% move"
Set-TextCell    $wsAll "M34" "time_exit_5min"
Set-NumCell     $wsAll "N34" 5

# Row 35 <- leadlag trade #34 (leadlag row 30)
Set-NumCell     $wsAll "A35" 34
Set-TextCell    $wsAll "B35" "2026-02-16"
Set-TextCell    $wsAll "C35" "21:28:34"
Set-TextCell    $wsAll "D35" "leadlag"
Set-TextCell    $wsAll "E35" "DOWN"
Set-NumCell     $wsAll "F35" 68810.35000000001
Set-NumCell     $wsAll "G35" 68703.52499799999
Set-TextCell    $wsAll "H35" "CLOSED"
Set-NumCell     $wsAll "I35" 0.1552
Set-NumCell     $wsAll "J35" 1.55
Set-NumCell     $wsAll "K35" 0.75
Set-TextCell    $wsAll "L35" "Binance leading with -0.125% move"
Set-TextCell    $wsAll "M35" "time_exit_5min"
Set-NumCell     $wsAll "N35" 5

# Row 36 <- momentum trade #35 (momentum row 7)
Set-NumCell     $wsAll "A36" 35
Set-TextCell    $wsAll "B36" "2026-02-16"
Set-TextCell    $wsAll "C36" "21:28:40"
Set-TextCell    $wsAll "D36" "momentum"
Set-TextCell    $wsAll "E36" "DOWN"
Set-NumCell     $wsAll "F36" 68718.965
Set-NumCell     $wsAll "G36" 69046.8011
Set-TextCell    $wsAll "H36" "CLOSED"
Set-NumCell     $wsAll "I36" -0.4771
Set-NumCell     $wsAll "J36" -4.77
Set-NumCell     $wsAll "K36" 0.9
Set-TextCell    $wsAll "L36" "Downward momentum: -0.299% over 10 samples"
Set-TextCell    $wsAll "M36" "time_exit_5min"
Set-NumCell     $wsAll "N36" 5

# Row 37 <- momentum trade #36 (momentum row 8)
Set-NumCell     $wsAll "A37" 36
Set-TextCell    $wsAll "B37" "2026-02-16"
Set-TextCell    $wsAll "C37" "21:28:47"
Set-TextCell    $wsAll "D37" "momentum"
Set-TextCell    $wsAll "E37" "DOWN"
Set-NumCell     $wsAll "F37" 68744.55
Set-NumCell     $wsAll "G37" 68269.16643899999
Set-TextCell    $wsAll "H37" "CLOSED"
Set-NumCell     $wsAll "I37" 0.6915
Set-NumCell     $wsAll "J37" 6.92
Set-NumCell     $wsAll "K37" 0.9
Set-TextCell    $wsAll "L37" "Downward momentum: -0.211% over 10 samples"
Set-TextCell    $wsAll "M37" "time_exit_5min"
Set-NumCell     $wsAll "N37" 5

# Row 38 <- momentum trade #37 (momentum row 9)
Set-NumCell     $wsAll "A38" 37
Set-TextCell    $wsAll "B38" "2026-02-16"
Set-TextCell    $wsAll "C38" "21:28:53"
Set-TextCell    $wsAll "D38" "momentum"
Set-TextCell    $wsAll "E38" "DOWN"
Set-NumCell     $wsAll "F38" 68582.45
Set-NumCell     $wsAll "G38" 68423.353225
Set-TextCell    $wsAll "H38" "CLOSED"
Set-NumCell     $wsAll "I38" 0.232
Set-NumCell     $wsAll "J38" 2.32
Set-NumCell     $wsAll "K38" 0.9
Set-TextCell    $wsAll "L38" "Downward momentum: -0.488% over 10 samples"
Set-TextCell    $wsAll "M38" "time_exit_5min"
Set-NumCell     $wsAll "N38" 5

# Row 39 <- momentum trade #38 (momentum row 10)
Set-NumCell     $wsAll "A39" 38
Set-TextCell    $wsAll "B39" "2026-02-16"
Set-TextCell    $wsAll "C39" "21:28:59"
Set-TextCell    $wsAll "D39" "momentum"
Set-TextCell    $wsAll "E39" "DOWN"
Set-NumCell     $wsAll "F39" 68584.08
Set-NumCell     $wsAll "G39" 68082.94107
Set-TextCell    $wsAll "H39" "CLOSED"
Set-NumCell     $wsAll "I39" 0.7307
Set-NumCell     $wsAll "J39" 7.31
Set-NumCell     $wsAll "K39" 0.9
Set-TextCell    $wsAll "L39" "Downward momentum: -0.524% over 10 samples"
Set-TextCell    $wsAll "M39" "time_exit_5min"
Set-NumCell     $wsAll "N39" 5

# ---------------------------------------------------------------------------
# 5) Comparison sheet
# ---------------------------------------------------------------------------
$wsComparison = $wb.Worksheets.Item("Comparison")

Set-NumCell  $wsComparison "B2" 47
Set-TextCell $wsComparison "C2" "36.2%"
Set-TextCell $wsComparison "D2" "2.19"
Set-TextCell $wsComparison "E2" "+0.5007%"
Set-TextCell $wsComparison "F2" "-0.3236%"
Set-TextCell $wsComparison "G2" "1.55"

Set-NumCell  $wsComparison "B3" 12
Set-TextCell $wsComparison "C3" "58.3%"
Set-TextCell $wsComparison "D3" "3.60"
Set-TextCell $wsComparison "E3" "+0.5778%"
Set-TextCell $wsComparison "F3" "-0.5622%"
Set-TextCell $wsComparison "G3" "1.03"
